$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the picture paths for the "Kinloch Heath" (row 6) and "Kinloch Luke" (row 5) dogs:
# strip the stray "Kinloch" segment that had crept into the file names.
# Set H6 before H5 so the shared-string table gets the new strings appended
# in the same order as the target workbook (Heath, then Luke).
$ws.Range("H6").Value = "PicturesOrig/FolkeNoertemann_Heath.jpg"
$ws.Range("H5").Value = "PicturesOrig/FolkeNoertemann_Luke.jpg"

# Move the selection/active cell to E13, matching the state the sheet was
# left in when the author saved the workbook.
$ws.Range("E13").Select()
